# Insert a new weekly price record at row 163, shifting all rows from the
# old row 163 onward down by one (old row 163 -> new row 164, ...,
# old row 235 -> new row 236). The new row 163 holds a brand-new entry.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("163:163").Insert()

$ws.Range("A163").Value = 3
$ws.Range("B163").Value = "Femacal de La Calera"
$ws.Range("C163").Value = "Coquimbo"
$ws.Range("D163").Value = 44489
$ws.Range("E163").Value = 5
$ws.Range("F163").Value = 100112003
$ws.Range("G163").Value = "Ajo"
$ws.Range("H163").Value = "Chino"
$ws.Range("I163").Value = "Primera"
$ws.Range("J163").Value = 73
$ws.Range("K163").Value = 16000
$ws.Range("L163").Value = 16500
$ws.Range("M163").Value = 16240
$ws.Range("N163").Value = "$/caja 10 kilos"
$ws.Range("O163").Value = "China"
$ws.Range("P163").Value = 1624
$ws.Range("Q163").Value = 10
$ws.Range("R163").Value = "Hortaliza"
